$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.988.06"
$ws.Range("E2").Value = "  -2.05%  "

$ws.Range("D3").Value = "2.169.61"
$ws.Range("E3").Value = "  -3.07%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.17"
$ws.Range("E5").Value = "  -1.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  -3.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.60"
$ws.Range("E7").Value = "  -7.56%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  -0.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.92"
$ws.Range("E10").Value = "  -0.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0925"
$ws.Range("E11").Value = "  -5.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.41"
$ws.Range("E12").Value = "  -16.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -2.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").Value = "2.493.53"
$ws.Range("E15").Value = "  -2.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.851"
$ws.Range("E16").Value = "  -1.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.21"
$ws.Range("E17").Value = "  -6.03%  "

$ws.Range("D18").Value = "2.207.21"
$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("D19").Value = "40.897.18"
$ws.Range("E19").Value = "  -2.11%  "

$ws.Range("E20").Value = "  -3.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.08"
$ws.Range("E21").Value = "  -2.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.49"
$ws.Range("E22").Value = "  -2.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.73"
$ws.Range("E23").Value = "  -2.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -10.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.56"
$ws.Range("E25").Value = "  +12.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.41"
$ws.Range("E28").Value = "  -4.22%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  -3.43%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.37"
$ws.Range("E30").Value = "  -2.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.20"
$ws.Range("E31").Value = "  -2.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -3.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.58"
$ws.Range("E33").Value = "  +2.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0735"
$ws.Range("E34").Value = "  +1.14%  "

$ws.Range("E35").Value = "  -3.60%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").Value = "  -3.66%  "

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.50"
$ws.Range("E37").Value = "  -4.03%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.05"
$ws.Range("E38").Value = "  -1.10%  "

$ws.Range("E39").Value = "  +6.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.16"
$ws.Range("E40").Value = "  -6.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.47"
$ws.Range("E41").Value = "  -10.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.49"
$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.23"
$ws.Range("E43").Value = "  -13.45%  "

$ws.Range("E44").Value = "  -7.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.191"
$ws.Range("E45").Value = "  -11.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.49"
$ws.Range("E46").Value = "  -5.38%  "

$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0986"
$ws.Range("E48").Value = "  -3.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.15"
$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  -3.81%  "

$ws.Range("E51").Value = "  -0.84%  "
